# Customize MCGLT to hit 1.5 target
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCGLT")

# Column C holds a per-technology annual capacity-growth cap (MW); column D
# mirrors it via a "=Cn" formula on most rows. Raise the cap from 500 to
# 1200 on every populated row.
$rows = @(2, 4, 6, 8, 10, 12, 14, 16, 18, 20, 22, 24, 28, 30, 32)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 1200
}

# Reposition the MCGLT window (scrolled down a bit further, new selection)
# and make MCGLT the active/selected sheet (it was previously "About").
$ws.Activate()
$ws.Range("C33").Select()
$win = $wb.Windows.Item(1)
$win.ScrollRow = 16
$win.ScrollColumn = 1
